# Renewable constraint feature: add per-project-site timeseries-column
# titles (title_demand / title_pv / title_wind) on the "project_sites"
# sheet, and point the example timeseries files at the relative
# ./inputs/ folder instead of the original author's local machine path.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project_sites")

# New header row (C2:E2) describing the extra timeseries columns.
$ws.Range("C2").Value = "title_demand"
$ws.Range("D2").Value = "title_pv"
$ws.Range("E2").Value = "title_wind"

# Masbate (row 3): fix timeseries path + new column headers to reference.
$ws.Range("B3").Value = "./inputs/Example1_Masbate.csv"
$ws.Range("C3").Value = "Demand"
$ws.Range("D3").Value = "SolarPV"
$ws.Range("E3").Value = "Wind"

# Atulayan (row 4): fix timeseries path + new column headers to reference.
$ws.Range("B4").Value = "./inputs/Example2_Atulayan.csv"
$ws.Range("C4").Value = "Demand"
$ws.Range("D4").Value = "SolarPV"
$ws.Range("E4").Value = "Wind"

# Restore the active selection on this sheet to B8.
[void]$ws.Range("B8").Select()
